$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting D:K to E:L (new FY2018 column)
$ws.Columns("D:D").Insert()

# Copy number formats/styles from the (now-shifted) old D column (now E) into new D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# Populate the new column D with FY2018 figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 1538600
$ws.Range("D9").Value = 1127900
$ws.Range("D10").Value = 410700
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 6700
$ws.Range("D15").Value = 4200
$ws.Range("D17").Value = 1431400
$ws.Range("D18").Value = 107200
$ws.Range("D20").Value = -6100
$ws.Range("D21").Value = 130300
$ws.Range("D22").Value = 21500
$ws.Range("D23").Value = 79600
$ws.Range("D24").Value = -2100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 81700
$ws.Range("D27").Value = 78900
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 2300
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 6100
$ws.Range("D33").Value = 81200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 81200
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 68800
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 379200
$ws.Range("D44").Value = 128800
$ws.Range("D45").Value = 21600
$ws.Range("D46").Value = 598400
$ws.Range("D47").Value = "NA"
$ws.Range("D48").Value = 184200
$ws.Range("D49").Value = 592800
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 682100
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2057500
$ws.Range("D57").Value = 153600
$ws.Range("D58").Value = 49900
$ws.Range("D59").Value = 266700
$ws.Range("D60").Value = 470200
$ws.Range("D61").Value = 331900
$ws.Range("D62").Value = 840500
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1642600
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -650100
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 414900
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 81200
$ws.Range("D83").Value = 29200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 110600
$ws.Range("D91").Value = -12400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -180600
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 16800
$ws.Range("D101").Value = -2300
$ws.Range("D102").Value = -55500
